# Apply the cryptos-list refresh: updates the Price (D) / Volume(1h) (E)
# columns for the refreshed rows, and fixes two row-ordering swaps
# (rows 12<->13 and rows 50<->51) by rewriting their Coin (B) / Link (C) /
# Price (D) / Volume(1h) (E) values.
#
# All of these source cells are plain text (inline strings) in the
# workbook, even the ones in the Price column that look like numbers
# (e.g. "1.00", "578.93"). Assigning such a string straight to .Value
# makes Excel infer a number, which both loses the intended formatting
# (trailing zero) and the original text type. For any new Price value
# that parses as a plain number we instead enter it with a leading
# apostrophe (Excel's standard "force text" quote-prefix) and then
# restore the cell's style to Normal so no stray number-format/quote-
# prefix styling is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.947.01'
$ws.Range("E2").Value = '  -1.63%  '
$ws.Range("D3").Value = '2.454.66'
$ws.Range("E3").Value = '  -2.78%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = '''578.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.88%  '
$ws.Range("D6").Value = '''165.67'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.12%  '
$ws.Range("D8").Value = '''0.512'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.27%  '
$ws.Range("D9").Value = '2.453.81'
$ws.Range("E9").Value = '  -2.81%  '
$ws.Range("E10").Value = '  -4.47%  '
$ws.Range("E11").Value = '  -1.17%  '
$ws.Range("B12").Value = 'Toncoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D12").Value = '''4.87'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.58%  '
$ws.Range("B13").Value = 'Cardano'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D13").Value = '''0.332'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.73%  '
$ws.Range("E14").Value = '  -4.76%  '
$ws.Range("D15").Value = '2.888.21'
$ws.Range("E15").Value = '  -3.28%  '
$ws.Range("D16").Value = '66.850.39'
$ws.Range("E16").Value = '  -1.72%  '
$ws.Range("E17").Value = '  -5.78%  '
$ws.Range("D18").Value = '2.447.93'
$ws.Range("E18").Value = '  -3.00%  '
$ws.Range("E19").Value = '  -5.14%  '
$ws.Range("E20").Value = '  -4.65%  '
$ws.Range("D21").Value = '''354.89'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.61%  '
$ws.Range("D22").Value = '''4.05'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.75%  '
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("D24").Value = '''69.49'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.50%  '
$ws.Range("D25").Value = '''4.22'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -9.50%  '
$ws.Range("D26").Value = '''1.76'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -8.18%  '
$ws.Range("D27").Value = '''8.93'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -11.06%  '
$ws.Range("D28").Value = '''0.997'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.34%  '
$ws.Range("D29").Value = '2.571.55'
$ws.Range("E29").Value = '  -2.99%  '
$ws.Range("E30").Value = '  -8.73%  '
$ws.Range("D31").Value = '''505.92'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.17%  '
$ws.Range("E32").Value = '  -6.70%  '
$ws.Range("E33").Value = '  -7.01%  '
$ws.Range("E34").Value = '  -8.19%  '
$ws.Range("D35").Value = '''1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("D36").Value = '''158.91'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.67%  '
$ws.Range("E37").Value = '  -9.22%  '
$ws.Range("D38").Value = '''18.49'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.48%  '
$ws.Range("E39").Value = '  -0.55%  '
$ws.Range("E40").Value = '  -7.09%  '
$ws.Range("E41").Value = '  +0.02%  '
$ws.Range("E42").Value = '  -7.47%  '
$ws.Range("E43").Value = '  -7.30%  '
$ws.Range("D44").Value = '''4.75'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -8.29%  '
$ws.Range("D45").Value = '''38.67'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.21%  '
$ws.Range("D46").Value = '''2.30'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -8.61%  '
$ws.Range("D47").Value = '''141.42'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.37%  '
$ws.Range("E48").Value = '  -6.36%  '
$ws.Range("E49").Value = '  -7.63%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '0.0₆0251'
$ws.Range("E50").Value = '  -10.11%  '
$ws.Range("B51").Value = 'Optimism'
$ws.Range("C51").Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range("D51").Value = '''1.58'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -8.68%  '
